$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> (new D price, new E volume) values to update, preserving text formatting
$updates = @{
    2 = @('310.54', '1.55%')
    3 = @('35.62', '-1.75%')
    4 = @('5.109', '1.07%')
    5 = @('0.08222', '4.55%')
    6 = @('2.062', '-4.69%')
    7 = @('7.939', '-0.59%')
    8 = @('2.962', '11.36%')
    9 = @('0.9268', '-0.14%')
    10 = @('0.1124', '15.04%')
    11 = @('0.1917', '2.62%')
    12 = @('0.09326', '4.01%')
    13 = @('0.03625', '-4.19%')
    14 = @('0.09916', '0.15%')
    15 = @('0.001446', '0.08%')
    16 = @('0.005828', '2.66%')
    17 = @($null, '0.04%')
    18 = @('4.129', '-0.67%')
    19 = @($null, '0.18%')
    20 = @('0.1310', '-0.31%')
    21 = @('5.095', '-0.85%')
    23 = @('0.04553', '-0.82%')
    24 = @('0.001226', '-0.55%')
    25 = @($null, '0.64%')
    26 = @('0.0001250', '-4.18%')
    27 = @('0.0004447', '-6.15%')
    39 = @('0.01980', '2.67%')
    40 = @('0.04923', '-1.67%')
    41 = @('0.007626', '-2.27%')
    42 = @('0.009958', '27.06%')
    43 = @($null, '-0.42%')
    44 = @('0.002129', '-0.88%')
    45 = @('0.01157', '2.85%')
    46 = @('0.00006548', '4.05%')
    47 = @($null, '-0.36%')
    48 = @('178.24', '244.37%')
    49 = @('0.001499', '-21.39%')
    50 = @('0.00002100', '-0.36%')
    51 = @('0.0002000', '-0.36%')
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $dVal = $vals[0]
    $eVal = $vals[1]
    if ($dVal -ne $null) {
        $dCell = $ws.Cells.Item($row, 4)
        $dCell.Value = "'" + $dVal
        $dCell.Style = "Normal"
    }
    $eCell = $ws.Cells.Item($row, 5)
    $eCell.Value = "'" + $eVal
    $eCell.Style = "Normal"
}
